$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: 2024-04-12_K_UHfit.dat ---
$ws.Range("B3").Value = "K"
$ws.Range("A3").Value = "2024-04-12_K_UHfit.dat"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 45394
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 1.8
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 202.14
$ws.Range("I3").Value = "2024-04-05_G"

# --- Row 4: 2024-04-12_L_UHfit.dat ---
$ws.Range("A4").Value = "2024-04-12_L_UHfit.dat"
$ws.Range("B4").Value = "L"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 45394
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 0.9
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 202.14
$ws.Range("I4").Value = "2024-03-05_E"

# --- Row 5: 2024-04-15_I_UHfit.dat ---
$ws.Range("B5").Value = "I"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 45397
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 0.9
$ws.Range("G5").Value = 10
$ws.Range("H5").Value = 202.14
$ws.Range("I5").Value = "2024-03-05_E"
$ws.Range("A5").Value = "2024-04-15_I_UHfit.dat"

# Copy the date number format from the existing date cell (D2) onto the
# new date cells so they reuse the same style record instead of creating
# a brand-new one.
$ws.Range("D2").Copy()
$ws.Range("D3:D5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection to match the saved workbook state.
$ws.Range("H11").Select()

# Restore the window vertical position recorded in the workbook view.
$excel.ActiveWindow.Top = 740
